$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "1.00" or "0.600"
# keep their literal formatting instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.457.89'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '2.606.55'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '537.37'
$ws.Range("E5").Value = '  +3.11%  '
$ws.Range("D6").Value = '141.48'
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").Value = '6.51'
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").Value = '59.385.72'
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").Value = '20.68'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '2.639.22'
$ws.Range("E16").Value = '  +2.01%  '
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = '341.36'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '67.49'
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +6.04%  '
$ws.Range("D31").Value = '5.81'
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").Value = '18.85'
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").Value = '149.83'
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").Value = '3.97'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = '0.833'
$ws.Range("E37").Value = '  +3.54%  '
$ws.Range("D38").Value = '0.826'
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '273.49'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").Value = '0.600'
$ws.Range("E42").Value = '  +2.36%  '
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").Value = '1.948.89'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").Value = '18.50'
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("D50").Value = '111.00'
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("E51").Value = '  +0.43%  '

# Restore the original (default/Normal) cell style on column D so only the
# text content changed, matching the source workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
